$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# Extend the training sheet with a new "TRIGGER" Q/A block (rows
# 172-174 and 176-179, with a blank separator row 175 as in the
# original pattern used elsewhere in the sheet, e.g. rows 111-116 /
# 117 blank).
#
# Cell-write order below matters: it controls the order in which new
# values are interned into xl/sharedStrings.xml, so the order here is
# deliberately NOT simple row-by-row — it reproduces the exact shared
# string allocation order of the target workbook.
# --------------------------------------------------------------------

# New, previously-unseen strings first, in the exact order they must
# be interned as new shared-string entries.
$ws.Range("C172").Value = "TRIGGER_FACEAPP_ON"
$ws.Range("C176").Value = "TRIGGER_FACEAPP_OFF"
$ws.Range("E172").Value = "TRIGGER"
$ws.Range("C177").Value = "Hey, wer hat das List aus gemacht?"
$ws.Range("C178").Value = "{name} warum haben Sie mein Display deaktiviert?"
$ws.Range("C179").Value = "Das ist aber nicht nett {name}. Machen Sie bitte mein Display wieder an."
$ws.Range("C173").Value = "Danke, dass Sie mein Display wieder eingeschatet haben."
$ws.Range("C174").Value = "Danke! Sehr nett von Ihnen."

# Remaining cells reuse already-existing shared strings ("Q", "A",
# "TRIGGER", "User:Unknown"); order no longer affects sharedStrings
# allocation.
$ws.Range("A172").Value = "Q"

$ws.Range("A173").Value = "A"
$ws.Range("D173").Value = "TRIGGER"

$ws.Range("A174").Value = "A"
$ws.Range("D174").Value = "TRIGGER"

$ws.Range("A176").Value = "Q"
$ws.Range("E176").Value = "TRIGGER"

$ws.Range("A177").Value = "A"
$ws.Range("B177").Value = "User:Unknown"
$ws.Range("D177").Value = "TRIGGER"

$ws.Range("A178").Value = "A"
$ws.Range("D178").Value = "TRIGGER"

$ws.Range("A179").Value = "A"
$ws.Range("D179").Value = "TRIGGER"

# --------------------------------------------------------------------
# Restore the view state (selection moves to the newly-entered block,
# scrolled so the new rows are visible).
# --------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 167
$win.ScrollColumn = 1
$ws.Range("C175").Select() | Out-Null
